$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sponsored`n1 BHK by Lodha® in Thane - 1,2,3 BHK by Lodha® in Thane`nLodha Group`nhttps://www.lodhagroup.in`nWorld-class 1 BHKs in Thane by India's #1 real estate developer. Building a better life. World-class homes in Thane by India's #1 real estate developer. Building a better life."
$ws.Range("B2").Value = "Macrotech Developers Limited"

$ws.Range("A3").Value = "Sponsored`n2 BHK Projects in Thane West | Starts at ₹93 Lacs* by Runwal`nlandsend.runwal.com`nhttp://landsend.runwal.com › projects › thane`nTake Advantage of the Umbrella Offer: 2 BHK Flats Starts at ₹93L* at Lands End by Runwal. Book Your Dream Home at Runwal Lands End And Avail Pay 10% Now & 90% On...`nView Location · View Gallery · Lands End by Runwal · Location Advantages"

$ws.Range("A4").Value = "Sponsored`nLaunching 2 BHK in Thane West | Pay 20% & Nothing till Jan'25`nraymondtenxera.com`nhttps://www.raymondtenxera.com › thane › project`nNew Launch Homes by Raymond with 38 Habitable Floors, 26500 SqFt Clubhouse, 40+ Amenities. Experience a futuristic lifestyle with Raymond Realty's Spacious 2 BHK homes in Thane West. 38 Storey Three towers. 26,500 sq.ft Clubhouse. Library & Reading Lounges.`nLocation Map · Overview · Location Advantages · Configuration · View Amenities · About Us"
$ws.Range("B4").Value = "Raymond Limited"

$ws.Range("A5").Value = "Sponsored`nGodrej Ascend, Kolshet, Thane | 2&3 BHK at ₹1.09Cr+*(All Incl)`nGodrej Properties`nhttps://www.godrejproperties.com`nReach Airoli in 20 min* and avail easy connectivity to Navi-Mumbai via Thane-Belapur Road. Book 2&3 BHK Opulent Residences at ₹1.09Cr+(All Incl)* with 40+ lifestyle amenities."

$ws.Range("A6").Value = "Sponsored`nFlats in Thane West | View Prices, Deals & Offers`nHousing.com`nhttps://www.housing.com › -- › --`nView listing photos, floor plan & use our detailed filters to find the perfect place. Looking for Property in Thane West? Housing.com offers 6768+ options in Thane West. Price History & Trends. Real & Verified Photos.`nLuxury Flats in Mumbai · Apartments In Thane West · Townships In Thane West"
$ws.Range("B6").Value = "Locon Solutions Pvt. Ltd."
